# Apply the "automatic update" edit:
#  1) Column C (the "Förändrad" / last-changed date) is bumped from
#     serial 45184 to 45186 for every data row (rows 2..201).
#  2) The HYPERLINK() formulas in columns S, T, U, V, W, X, Y (present only
#     on rows 2..11) get a friendly-text second argument equal to the
#     row's "Beteckning" value (column A), e.g.
#       =HYPERLINK("...A 57042-2021.xlsx")
#     becomes
#       =HYPERLINK("...A 57042-2021.xlsx", "A 57042-2021")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 201
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    # 1) Bump the "changed" date in column C.
    $ws.Cells.Item($r, 3).Value = 45186

    # 2) Add the display-text argument to any HYPERLINK formulas on this row.
    $label = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $r)
        $formula = $cell.Formula

        if ($formula -ne "" -and $formula.ToUpper().Contains("HYPERLINK(") -and -not $formula.Contains(",")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
